# Roll the quarterly income statement forward by one quarter:
#   - drop the oldest quarter (Q2 1399/06)
#   - shift every remaining quarter one column to the left (D<-E, E<-F, ... L<-M)
#   - populate the newly freed last column (M) with the new quarter
#     (Q4 1401/12, published 1402-02-27) and its figures
#   - a few historical figures were also restated when the "read_price"
#     algorithm changed (column I gets a corrected value instead of a
#     pure left-shift in rows 11, 12, 13 and 16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(4,5,6,7,8,9,10,11,12,13)   # D..M

# Row 8: quarter/period labels
$periods = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(8, $cols[$i]).Value = $periods[$i]
}

# Row 9: publish-date labels.
# Leading apostrophe forces these to stay text - without it, bare
# "yyyy-mm-dd"-shaped strings (no trailing "(n)") get auto-coerced into
# date serials by the value-type inference instead of staying strings.
$dates = @(
    "'1400-10-29 (2)",
    "'1401-03-24 (9)",
    "'1401-04-30 (3)",
    "'1401-08-29 (4)",
    "'1401-10-28 (2)",
    "'1402-02-27 (7)",
    "'1401-04-30",
    "'1401-08-29 (2)",
    "'1401-10-28",
    "'1402-02-27"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(9, $cols[$i]).Value = $dates[$i]
}

# Numeric data rows 11-27 (row 10 is a blank separator, left untouched)
$rowData = @{
    11 = @(282124, 273893, 217182, 332113, 470167, 612809, 590461, 779614, 760968, 868695)
    12 = @(-198794, -189752, -139110, -233480, -346561, -420397, -460541, -649763, -636680, -640754)
    13 = @(83330, 84141, 78072, 98633, 123606, 192412, 129920, 129851, 124288, 227941)
    14 = @(-7546, -8906, -14060, -16456, -11974, -18222, -14340, -23010, -23176, -34827)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(3361, 115, 999, 3173, -1900, 1657, 2795, 3329, 2664, 25271)
    17 = @(79145, 75350, 65011, 85350, 109732, 175847, 118375, 110170, 103776, 218385)
    18 = @(-13134, -16598, -28433, -21636, -35608, -83174, -30208, -46582, -53823, -132312)
    19 = @(3254, -2326, 671, -1867, -33289, 49468, 6848, 7418, 9337, 5749)
    20 = @(69265, 56426, 37249, 61847, 40835, 142141, 95015, 71006, 59290, 91822)
    21 = @(-11538, -10282, -6401, -12915, -4139, -33584, -13575, -18861, -1014, -6602)
    22 = @(57727, 46144, 30848, 48932, 36696, 108557, 81440, 52145, 58276, 85220)
    23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    24 = @(57727, 46144, 30848, 48932, 36696, 108557, 81440, 52145, 58276, 85220)
    25 = @(74, 60, 40, 63, 47, 140, 105, 67, 75, 110)
    26 = @(775000, 775000, 775000, 775000, 775000, 775000, 775000, 775000, 775000, 775000)
    27 = @(74, 60, 40, 63, 47, 140, 105, 67, 75, 110)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}

Write-Host "Quarterly roll-forward applied"
